$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), styled like the other header cells (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New data for columns I (I0) and J (IF), rows 2-43
$values = @(
    @(2, 5, 6),
    @(3, 6, 6),
    @(4, 4, 5),
    @(5, 7, 7),
    @(6, 8, 8),
    @(7, 6, 6),
    @(8, 7, 7),
    @(9, 7, 7),
    @(10, 8, 8),
    @(11, 5, 6),
    @(12, 7, 7),
    @(13, 11, 11),
    @(14, 3, 4),
    @(15, 10, 10),
    @(16, 7, 8),
    @(17, 11, 14),
    @(18, 7, 8),
    @(19, 7, 8),
    @(20, 7, 7),
    @(21, 7, 7),
    @(22, 4, 5),
    @(23, 6, 6),
    @(24, 6, 6),
    @(25, 5, 6),
    @(26, 6, 6),
    @(27, 8, 8),
    @(28, 8, 8),
    @(29, 3, 5),
    @(30, 8, 8),
    @(31, 8, 8),
    @(32, 9, 9),
    @(33, 6, 6),
    @(34, 9, 9),
    @(35, 8, 9),
    @(36, 7, 7),
    @(37, 5, 5),
    @(38, 8, 8),
    @(39, 9, 9),
    @(40, 6, 6),
    @(41, 6, 6),
    @(42, 8, 8),
    @(43, 4, 4)
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}

Write-Output "Added I0/IF columns"
